# Team Attendance update (b8-B1, 09/09/23) — add three new daily rows
# (07-Sep-2023, 08-Sep-2023, 09-Sep-2023) with PRESENT/ABSENT marks and
# the "no response / out of town / outside / stuck in rain / ..." remark
# comments that go with the new ABSENT cells.

# Comments in this workbook are authored by "LENOVO" - match that so any
# newly added notes line up with the existing ones.
$excel.UserName = "LENOVO"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 33: 07-Sep-2023 ------------------------------------------------
$ws.Range("A33").NumberFormat = "d-mmm-yy"
$ws.Range("A33").Value = (Get-Date -Year 2023 -Month 9 -Day 7 -Hour 0 -Minute 0 -Second 0)
$ws.Range("B33:F33").Value = "PRESENT"
$ws.Range("G33:I33").Value = "ABSENT"
$ws.Range("J33:K33").Value = "PRESENT"

[void]$ws.Range("G33").AddComment("LENOVO:" + [char]10 + "no reponse")
[void]$ws.Range("I33").AddComment("LENOVO:" + [char]10 + "out of town" + [char]10)

# --- Row 34: 08-Sep-2023 -------------------------------------------------
$ws.Range("A34").NumberFormat = "d-mmm-yy"
$ws.Range("A34").Value = (Get-Date -Year 2023 -Month 9 -Day 8 -Hour 0 -Minute 0 -Second 0)
$ws.Range("B34").Value = "PRESENT"
$ws.Range("C34:D34").Value = "ABSENT"
$ws.Range("E34:F34").Value = "PRESENT"
$ws.Range("G34").Value = "PRESENT"
$ws.Range("H34:K34").Value = "ABSENT"

[void]$ws.Range("C34").AddComment("LENOVO:" + [char]10 + "outide")
[void]$ws.Range("H34").AddComment("LENOVO:" + [char]10 + "outsside")
[void]$ws.Range("J34").AddComment("LENOVO:" + [char]10 + "Stuck in rain")
[void]$ws.Range("K34").AddComment("LENOVO:" + [char]10 + "Some work")

# --- Row 35: 09-Sep-2023 -------------------------------------------------
$ws.Range("A35").NumberFormat = "d-mmm-yy"
$ws.Range("A35").Value = (Get-Date -Year 2023 -Month 9 -Day 9 -Hour 0 -Minute 0 -Second 0)
$ws.Range("B35:F35").Value = "PRESENT"
$ws.Range("G35:I35").Value = "ABSENT"
$ws.Range("J35").Value = "PRESENT"
$ws.Range("K35").Value = "ABSENT"

# Leave the same cell selected that Excel would land on after typing the
# last row of the table.
[void]$ws.Range("K35").Select()
